# Refresh the cryptos list snapshot scraped by the GitHub Actions job:
# updates the Price / Volume(1h) columns for the existing rows and applies
# the reshuffle at the bottom of the table (row 49: Algorand -> BabyDogeCoin,
# row 51: EnergySwap -> the refreshed Algorand entry).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D stores the price as a literal display string ("27.101.48",
# "19.50", ...). Plain-decimal-looking values (e.g. "19.50", "1.00") would
# otherwise be auto-parsed into numbers by Excel and lose their formatting,
# so those specific cells are pre-formatted as Text before the write.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

$ws.Range('D2').Value = '27.101.48'
$ws.Range('E2').Value = '  +3.15%  '
$ws.Range('D3').Value = '1.655.49'
$ws.Range('E3').Value = '  +3.51%  '
Set-TextValue 'D4' '0.999'
$ws.Range('E4').Value = '  -0.16%  '
Set-TextValue 'D5' '215.29'
$ws.Range('E5').Value = '  +1.23%  '
Set-TextValue 'D6' '0.509'
$ws.Range('E6').Value = '  +1.32%  '
Set-TextValue 'D7' '0.999'
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('E8').Value = '  +2.10%  '
$ws.Range('E9').Value = '  +1.33%  '
Set-TextValue 'D10' '19.50'
$ws.Range('E11').Value = '  +0.54%  '
$ws.Range('D12').Value = '1.888.99'
$ws.Range('E12').Value = '  +3.47%  '
$ws.Range('D13').Value = '1.657.80'
$ws.Range('E13').Value = '  +2.49%  '
$ws.Range('E14').Value = '  +1.61%  '
Set-TextValue 'D15' '0.520'
$ws.Range('E15').Value = '  +2.76%  '
Set-TextValue 'D16' '64.95'
$ws.Range('E16').Value = '  +1.97%  '
Set-TextValue 'D17' '240.37'
$ws.Range('E17').Value = '  +4.98%  '
$ws.Range('D18').Value = '27.071.41'
$ws.Range('E18').Value = '  +3.01%  '
Set-TextValue 'D19' '7.84'
$ws.Range('E19').Value = '  +3.00%  '
$ws.Range('D20').Value = '0.0₃0728'
$ws.Range('E20').Value = '  +1.13%  '
Set-TextValue 'D21' '0.998'
$ws.Range('E21').Value = '  -0.18%  '
Set-TextValue 'D22' '4.46'
$ws.Range('E22').Value = '  +4.20%  '
Set-TextValue 'D24' '9.27'
$ws.Range('E24').Value = '  +3.53%  '
Set-TextValue 'D25' '146.10'
$ws.Range('E25').Value = '  +0.33%  '
Set-TextValue 'D26' '1.00'
$ws.Range('E26').Value = '  -0.16%  '
Set-TextValue 'D27' '7.11'
$ws.Range('E27').Value = '  +2.49%  '
$ws.Range('E28').Value = '  +1.26%  '
Set-TextValue 'D29' '15.84'
$ws.Range('E29').Value = '  +2.78%  '
$ws.Range('E30').Value = '  +0.48%  '
$ws.Range('E31').Value = '  +0.59%  '
$ws.Range('D32').Value = '1.520.97'
$ws.Range('E32').Value = '  +5.30%  '
$ws.Range('E33').Value = '  +2.72%  '
Set-TextValue 'D34' '3.04'
$ws.Range('E34').Value = '  +2.78%  '
$ws.Range('E35').Value = '  +6.54%  '
$ws.Range('E36').Value = '  -0.20%  '
Set-TextValue 'D37' '0.577'
$ws.Range('E37').Value = '  +1.38%  '
Set-TextValue 'D38' '0.894'
$ws.Range('E38').Value = '  +8.86%  '
$ws.Range('E39').Value = '  +2.46%  '
$ws.Range('E40').Value = '  +3.05%  '
Set-TextValue 'D41' '0.999'
$ws.Range('E41').Value = '  -0.25%  '
$ws.Range('E42').Value = '  +4.22%  '
Set-TextValue 'D43' '64.69'
$ws.Range('D44').Value = '1.795.38'
$ws.Range('E44').Value = '  +3.21%  '
Set-TextValue 'D45' '0.771'
Set-TextValue 'D46' '0.918'
$ws.Range('E46').Value = '  -0.61%  '
Set-TextValue 'D47' '90.49'
$ws.Range('E47').Value = '  +3.47%  '
$ws.Range('E48').Value = '  +3.14%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0103'
$ws.Range('E49').Value = '  -2.26%  '
$ws.Range('E50').Value = '  +0.67%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D51' '0.0978'
$ws.Range('E51').Value = '  +2.86%  '
